$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells are treated as text so values like "1.00" or "0.0000282"
# are preserved exactly instead of being coerced into numbers.
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '65.678.95'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +1.34%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.399.53'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +0.71%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '560.63'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +0.48%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '176.26'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  -0.48%  '
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +1.07%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '3.389.78'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +0.67%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '1.00'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +0.01%  '
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +4.12%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.639'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +0.92%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '54.00'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  -2.18%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000282'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +1.26%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '9.24'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +1.10%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.934.51'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +0.72%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '18.42'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +0.62%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '3.383.86'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +0.38%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '65.556.72'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +1.42%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '11.91'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +0.13%  '
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +1.50%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '460.85'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +0.64%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '4.92'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +2.25%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '4.14'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +0.97%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '14.12'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +4.90%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '87.58'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +1.73%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.93'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  +2.39%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '10.74'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  -1.62%  '
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -0.76%  '
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +3.54%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '6.59'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -2.41%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '63.56'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +6.71%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '11.54'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  +0.19%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '581.93'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  -0.11%  '
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -0.17%  '
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -0.01%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '3.61'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +3.58%  '
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +0.65%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '36.09'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +0.15%  '
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +0.48%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0₃0748'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  -2.36%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '3.105.32'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  +0.04%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0420'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +1.67%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.80'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -2.23%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.46'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  -2.44%  '
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  +1.94%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.18'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -0.89%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.998'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +0.02%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '140.68'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  +3.74%  '
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -2.33%  '
$ws.Range('B51').NumberFormat = "@"
$ws.Range('B51').Value = 'THORChain'
$ws.Range('C51').NumberFormat = "@"
$ws.Range('C51').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '8.41'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +1.04%  '
